$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.385.60'
$ws.Range('D3').Value = '2.107.69'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.98'
$ws.Range('E5').Value = '  +2.60%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5231'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4568'
$ws.Range('E8').Value = '  +6.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.91'
$ws.Range('E9').Value = '  +16.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08919'
$ws.Range('E10').Value = '  +2.85%  '
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.44'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').Value = '2.096.86'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('E14').Value = '  +3.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.987'
$ws.Range('E15').Value = '  +5.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.52'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001135'
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06628'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.27'
$ws.Range('E20').Value = '  +3.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.368'
$ws.Range('E22').Value = '  +3.35%  '
$ws.Range('D23').Value = '30.447.68'
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.42'
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('E25').Value = '  +4.13%  '
$ws.Range('D26').Value = '2.349.24'
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.36'
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.571'
$ws.Range('E28').Value = '  +3.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.55'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.68'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.243'
$ws.Range('E31').Value = '  +6.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.723'
$ws.Range('E32').Value = '  +16.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1073'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E34').Value = '  +4.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.926'
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.53'
$ws.Range('E36').Value = '  +10.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02587'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06828'
$ws.Range('E38').Value = '  +4.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.557'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.77'
$ws.Range('E40').Value = '  +3.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2298'
$ws.Range('E41').Value = '  +3.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6907'
$ws.Range('E42').Value = '  +4.75%  '
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.352'
$ws.Range('E44').Value = '  +8.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.00'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6381'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.657'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000353'
$ws.Range('E49').Value = '  +26.52%  '
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3441'
$ws.Range('E51').Value = '  +28.45%  '
